# Update line-power-flow results (pl_mw) on Sheet1 for the 380 kV case.
# Columns B,D,E,F,G,I,J,L,N (rows 2-25) get new values; A,C,H,K,M,O are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.313917936120617
$ws.Range("D2").Value = 0.2620200641667338
$ws.Range("E2").Value = 0.1774925117578334
$ws.Range("F2").Value = 4.971185419538699
$ws.Range("G2").Value = 0.002632765969952425
$ws.Range("I2").Value = 1.490529160807945
$ws.Range("J2").Value = 0.1788591254052818
$ws.Range("L2").Value = 2.014708627818834
$ws.Range("N2").Value = 1.766612247303891

$ws.Range("B3").Value = 2.222570831146925
$ws.Range("D3").Value = 0.2336770182094767
$ws.Range("E3").Value = 0.1548528657272996
$ws.Range("F3").Value = 4.905277031721369
$ws.Range("G3").Value = 0.002643053144961272
$ws.Range("I3").Value = 1.504587246263817
$ws.Range("J3").Value = 0.1553932956667694
$ws.Range("L3").Value = 1.924069595112741
$ws.Range("N3").Value = 1.792117770703404

$ws.Range("B4").Value = 2.167718327581611
$ws.Range("D4").Value = 0.216451060731373
$ws.Range("E4").Value = 0.1410105468693033
$ws.Range("F4").Value = 4.868492898754369
$ws.Range("G4").Value = 0.002649690116448375
$ws.Range("I4").Value = 1.514126605686357
$ws.Range("J4").Value = 0.1409970917365939
$ws.Range("L4").Value = 1.869725732864708
$ws.Range("N4").Value = 1.808506068148271

$ws.Range("B5").Value = 2.145673982405469
$ws.Range("D5").Value = 0.2094736778140884
$ws.Range("E5").Value = 0.1353830910397207
$ws.Range("F5").Value = 4.854420190662495
$ws.Range("G5").Value = 0.002652475716051317
$ws.Range("I5").Value = 1.518241123011919
$ws.Range("J5").Value = 0.1351324067345274
$ws.Range("L5").Value = 1.847904700745204
$ws.Range("N5").Value = 1.815367360036326

$ws.Range("B6").Value = 2.142032109718059
$ws.Range("D6").Value = 0.2083175843448544
$ws.Range("E6").Value = 0.1344494278149142
$ws.Range("F6").Value = 4.852138554402643
$ws.Range("G6").Value = 0.002652943164440202
$ws.Range("I6").Value = 1.518938023580191
$ws.Range("J6").Value = 0.1341586621069268
$ws.Range("L6").Value = 1.844300803942076
$ws.Range("N6").Value = 1.816517718480336

$ws.Range("B7").Value = 2.16741978339212
$ws.Range("D7").Value = 0.2163567926687904
$ws.Range("E7").Value = 0.1409346004871779
$ws.Range("F7").Value = 4.868299407910115
$ws.Range("G7").Value = 0.002649727355787472
$ws.Range("I7").Value = 1.514181177188391
$ws.Range("J7").Value = 0.1409179923622474
$ws.Range("L7").Value = 1.869430138754296
$ws.Range("N7").Value = 1.808597861358169

$ws.Range("B8").Value = 2.28216384015019
$ws.Range("D8").Value = 0.252209394514523
$ws.Range("E8").Value = 0.1696734131346815
$ws.Range("F8").Value = 4.947689980905352
$ws.Range("G8").Value = 0.002636246664829915
$ws.Range("I8").Value = 1.495187351301681
$ws.Range("J8").Value = 0.1707649539434044
$ws.Range("L8").Value = 1.983181803761511
$ws.Range("N8").Value = 1.775255445352125

$ws.Range("B9").Value = 2.517085914787799
$ws.Range("D9").Value = 0.324019362812038
$ws.Range("E9").Value = 0.2265572093989192
$ws.Range("F9").Value = 5.133042437469499
$ws.Range("G9").Value = 0.002612338376443874
$ws.Range("I9").Value = 1.465191974681971
$ws.Range("J9").Value = 0.2294428915165838
$ws.Range("L9").Value = 2.216858979302856
$ws.Range("N9").Value = 1.715648948492227

$ws.Range("B10").Value = 2.695902656930059
$ws.Range("D10").Value = 0.3778405792748458
$ws.Range("E10").Value = 0.2687632546904695
$ws.Range("F10").Value = 5.287946111611348
$ws.Range("G10").Value = 0.002596290393978269
$ws.Range("I10").Value = 1.447645423281699
$ws.Range("J10").Value = 0.2727234025321366
$ws.Range("L10").Value = 2.395348429833405
$ws.Range("N10").Value = 1.67538168529514

$ws.Range("B11").Value = 2.778643234879098
$ws.Range("D11").Value = 0.4025886952803432
$ws.Range("E11").Value = 0.288074415981086
$ws.Range("F11").Value = 5.362625697516904
$ws.Range("G11").Value = 0.002589314216436866
$ws.Range("I11").Value = 1.440653970721286
$ws.Range("J11").Value = 0.2924676650524987
$ws.Range("L11").Value = 2.478102931139745
$ws.Range("N11").Value = 1.657829552700539

$ws.Range("B12").Value = 2.810178796560933
$ws.Range("D12").Value = 0.4120008609106094
$ws.Range("E12").Value = 0.2954047381571883
$ws.Range("F12").Value = 5.391522498863083
$ws.Range("G12").Value = 0.002586718734637649
$ws.Range("I12").Value = 1.438150273054184
$ws.Range("J12").Value = 0.2999537493398066
$ws.Range("L12").Value = 2.50967016559207
$ws.Range("N12").Value = 1.651293275209568

$ws.Range("B13").Value = 2.803377953028928
$ws.Range("D13").Value = 0.4099719360816607
$ws.Range("E13").Value = 0.2938252134849648
$ws.Range("F13").Value = 5.385271417479117
$ws.Range("G13").Value = 0.00258727566704309
$ws.Range("I13").Value = 1.438683070772065
$ws.Range("J13").Value = 0.2983410489863445
$ws.Range("L13").Value = 2.502861278330556
$ws.Range("N13").Value = 1.652696068458443

$ws.Range("B14").Value = 2.781233589796727
$ws.Range("D14").Value = 0.4033622114791058
$ws.Range("E14").Value = 0.2886771239531782
$ws.Range("F14").Value = 5.364990617353016
$ws.Range("G14").Value = 0.002589099759955016
$ws.Range("I14").Value = 1.440445099789436
$ws.Range("J14").Value = 0.2930833542633309
$ws.Range("L14").Value = 2.480695345341815
$ws.Range("N14").Value = 1.657289596557128

$ws.Range("B15").Value = 2.767696113833154
$ws.Range("D15").Value = 0.3993189271794222
$ws.Range("E15").Value = 0.2855261108150984
$ws.Range("F15").Value = 5.352648771882798
$ws.Range("G15").Value = 0.002590223080990742
$ws.Range("I15").Value = 1.441543163536032
$ws.Range("J15").Value = 0.2898641240594202
$ws.Range("L15").Value = 2.467148185172391
$ws.Range("N15").Value = 1.660117641168972

$ws.Range("B16").Value = 2.690523689177098
$ws.Range("D16").Value = 0.3762287717142385
$ws.Range("E16").Value = 0.2675036012131073
$ws.Range("F16").Value = 5.28315137298344
$ws.Range("G16").Value = 0.002596752794709619
$ws.Range("I16").Value = 1.448122380286726
$ws.Range("J16").Value = 0.2714343042763119
$ws.Range("L16").Value = 2.389972082320639
$ws.Range("N16").Value = 1.676544184940287

$ws.Range("B17").Value = 2.643540370632536
$ws.Range("D17").Value = 0.3621332636390378
$ws.Range("E17").Value = 0.2564770248017822
$ws.Range("F17").Value = 5.241604079198595
$ws.Range("G17").Value = 0.002600841321238635
$ws.Range("I17").Value = 1.452413188661296
$ws.Range("J17").Value = 0.2601434280354056
$ws.Range("L17").Value = 2.343030313227814
$ws.Range("N17").Value = 1.686817600200447

$ws.Range("B18").Value = 2.616648028721443
$ws.Range("D18").Value = 0.3540506470038167
$ws.Range("E18").Value = 0.2501452197600855
$ws.Range("F18").Value = 5.218103032366173
$ws.Range("G18").Value = 0.002603223464365556
$ws.Range("I18").Value = 1.454974300408963
$ws.Range("J18").Value = 0.2536543784349305
$ws.Range("L18").Value = 2.316177011388561
$ws.Range("N18").Value = 1.692798648260933

$ws.Range("B19").Value = 2.60756519870688
$ws.Range("D19").Value = 0.3513181857036898
$ws.Range("E19").Value = 0.2480031163775607
$ws.Range("F19").Value = 5.210213627710175
$ws.Range("G19").Value = 0.002604035271881482
$ws.Range("I19").Value = 1.455857406549903
$ws.Range("J19").Value = 0.2514581442915755
$ws.Range("L19").Value = 2.307109910311681
$ws.Range("N19").Value = 1.694836096323913

$ws.Range("B20").Value = 2.648528224659628
$ws.Range("D20").Value = 0.3636311751065193
$ws.Range("E20").Value = 0.2576497360266785
$ws.Range("F20").Value = 5.245985808077052
$ws.Range("G20").Value = 0.002600402933472026
$ws.Range("I20").Value = 1.451946774055195
$ws.Range("J20").Value = 0.2613448170356492
$ws.Range("L20").Value = 2.348012155666424
$ws.Range("N20").Value = 1.685716518922709

$ws.Range("B21").Value = 2.787732380468583
$ws.Range("D21").Value = 0.4053025255672367
$ws.Range("E21").Value = 0.2901887520305451
$ws.Range("F21").Value = 5.370930735017708
$ws.Range("G21").Value = 0.002588562727464761
$ws.Range("I21").Value = 1.439923634979948
$ws.Range("J21").Value = 0.294627401692992
$ws.Range("L21").Value = 2.487199730947225
$ws.Range("N21").Value = 1.655937369609092

$ws.Range("B22").Value = 2.879897871494563
$ws.Range("D22").Value = 0.4327749881819329
$ws.Range("E22").Value = 0.3115582586856931
$ws.Range("F22").Value = 5.456192877705462
$ws.Range("G22").Value = 0.002581093849436667
$ws.Range("I22").Value = 1.432904736667339
$ws.Range("J22").Value = 0.3164347082991696
$ws.Range("L22").Value = 2.579509699900598
$ws.Range("N22").Value = 1.637118206935824

$ws.Range("B23").Value = 2.830597809780159
$ws.Range("D23").Value = 0.4180897981243845
$ws.Range("E23").Value = 0.3001429572384779
$ws.Range("F23").Value = 5.410353298413185
$ws.Range("G23").Value = 0.002585055602818863
$ws.Range("I23").Value = 1.436573631219446
$ws.Range("J23").Value = 0.304790231449573
$ws.Range("L23").Value = 2.530117312120126
$ws.Range("N23").Value = 1.647103410445236

$ws.Range("B24").Value = 2.646272845421947
$ws.Range("D24").Value = 0.3629539038679184
$ws.Range("E24").Value = 0.2571195299759665
$ws.Range("F24").Value = 5.244003629979545
$ws.Range("G24").Value = 0.002600601030413184
$ws.Range("I24").Value = 1.452157346575632
$ws.Range("J24").Value = 0.2608016622376681
$ws.Range("L24").Value = 2.345759447296473
$ws.Range("N24").Value = 1.686214085035116

$ws.Range("B25").Value = 2.452454350248559
$ws.Range("D25").Value = 0.3044170520095406
$ws.Range("E25").Value = 0.2111028156518699
$ws.Range("F25").Value = 5.079658507099396
$ws.Range("G25").Value = 0.00261853803714253
$ws.Range("I25").Value = 1.472523037621826
$ws.Range("J25").Value = 0.2135450002162003
$ws.Range("L25").Value = 2.152472797736664
$ws.Range("N25").Value = 1.731155242505192
